# -----------------------------------------------------------------------
# Add a new "2022-Q4" quarter to the 002653-海思科 holdings workbook:
#   * insert a new "2022-Q4" sheet (fund-holdings detail) right after the
#     existing "2022-Q3" sheet's position, pushing every later quarter
#     sheet one slot further right;
#   * add the matching summary row at the top of the data in "总计".
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($Worksheet, [int]$Row, [int]$Col, [string]$Text)
    $cell = $Worksheet.Cells.Item($Row, $Col)
    # Force a genuine "Text" cell so number-looking strings (fund codes
    # with leading zeros, "12.47", "0.5225", ...) keep every digit instead
    # of being silently re-interpreted as numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

# ==========================================================================
# 1) "总计" sheet: insert the 2022-Q4 summary as the new row 2 and push the
#    previously-existing rows down by one.
# ==========================================================================

$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Cells.Item(7, 1).Value = 5
$totalSheet.Cells.Item(7, 2).Value = "2021-Q1"
$totalSheet.Cells.Item(7, 3).Value = 4
$totalSheet.Cells.Item(7, 4).Value = 0.05
$totalSheet.Cells.Item(6, 1).Copy()
$totalSheet.Cells.Item(7, 1).PasteSpecial(-4122)

$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(6, 2).Value = "2021-Q2"
$totalSheet.Cells.Item(6, 3).Value = 4
$totalSheet.Cells.Item(6, 4).Value = 0.08

$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(5, 3).Value = 5
$totalSheet.Cells.Item(5, 4).Value = 0.55

$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(4, 3).Value = 2
$totalSheet.Cells.Item(4, 4).Value = 0.09

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(3, 3).Value = 11
$totalSheet.Cells.Item(3, 4).Value = 0.93

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 8
$totalSheet.Cells.Item(2, 4).Value = 0.91

# ==========================================================================
# 2) Duplicate the "2022-Q3" detail sheet (so the new sheet inherits the
#    exact same header/row styling) immediately before itself, rename the
#    duplicate to "2022-Q4", trim it down to the 2022-Q4 row count, and
#    overwrite every cell with the 2022-Q4 fund-holdings data.
# ==========================================================================

$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($q3Sheet)
$q4Sheet = $wb.Worksheets.Item("2022-Q3 (2)")
$q4Sheet.Name = "2022-Q4"

# The source sheet has 11 data rows (rows 2-12); 2022-Q4 only has 8 (rows
# 2-9) so drop the now-unused tail rows.
$q4Sheet.Range("A10:H12").Delete()

# ---- row 2 ----
$q4Sheet.Cells.Item(2, 1).Value = 0
Set-TextValue $q4Sheet 2 2 "009693"
$q4Sheet.Cells.Item(2, 3).Value = "富国积极成长一年定期开放混合"
Set-TextValue $q4Sheet 2 4 "12.47"
Set-TextValue $q4Sheet 2 5 "98.05"
Set-TextValue $q4Sheet 2 6 "4.19"
Set-TextValue $q4Sheet 2 7 "0.5225"
$q4Sheet.Cells.Item(2, 8).Value = 5

# ---- row 3 ----
$q4Sheet.Cells.Item(3, 1).Value = 1
Set-TextValue $q4Sheet 3 2 "002300"
$q4Sheet.Cells.Item(3, 3).Value = "长盛医疗行业量化配置股票"
Set-TextValue $q4Sheet 3 4 "2.79"
Set-TextValue $q4Sheet 3 5 "92.63"
Set-TextValue $q4Sheet 3 6 "5.90"
Set-TextValue $q4Sheet 3 7 "0.1646"
$q4Sheet.Cells.Item(3, 8).Value = 7

# ---- row 4 ----
$q4Sheet.Cells.Item(4, 1).Value = 2
Set-TextValue $q4Sheet 4 2 "000684"
$q4Sheet.Cells.Item(4, 3).Value = "长盛养老健康产业灵活配置混合"
Set-TextValue $q4Sheet 4 4 "1.47"
Set-TextValue $q4Sheet 4 5 "88.14"
Set-TextValue $q4Sheet 4 6 "4.95"
Set-TextValue $q4Sheet 4 7 "0.0728"
$q4Sheet.Cells.Item(4, 8).Value = 6

# ---- row 5 ----
$q4Sheet.Cells.Item(5, 1).Value = 3
Set-TextValue $q4Sheet 5 2 "008412"
$q4Sheet.Cells.Item(5, 3).Value = "长盛竞争优势股票A"
Set-TextValue $q4Sheet 5 4 "0.83"
Set-TextValue $q4Sheet 5 5 "88.91"
Set-TextValue $q4Sheet 5 6 "5.28"
Set-TextValue $q4Sheet 5 7 "0.0438"
$q4Sheet.Cells.Item(5, 8).Value = 6

# ---- row 6 ----
$q4Sheet.Cells.Item(6, 1).Value = 4
Set-TextValue $q4Sheet 6 2 "005970"
$q4Sheet.Cells.Item(6, 3).Value = "国泰消费优选股票"
Set-TextValue $q4Sheet 6 4 "0.98"
Set-TextValue $q4Sheet 6 5 "93.52"
Set-TextValue $q4Sheet 6 6 "3.65"
Set-TextValue $q4Sheet 6 7 "0.0358"
$q4Sheet.Cells.Item(6, 8).Value = 9

# ---- row 7 ----
$q4Sheet.Cells.Item(7, 1).Value = 5
Set-TextValue $q4Sheet 7 2 "004945"
$q4Sheet.Cells.Item(7, 3).Value = "长信中证500指数增强A"
Set-TextValue $q4Sheet 7 4 "2.08"
Set-TextValue $q4Sheet 7 5 "92.79"
Set-TextValue $q4Sheet 7 6 "1.61"
Set-TextValue $q4Sheet 7 7 "0.0335"
$q4Sheet.Cells.Item(7, 8).Value = 6

# ---- row 8 ----
$q4Sheet.Cells.Item(8, 1).Value = 6
Set-TextValue $q4Sheet 8 2 "008413"
$q4Sheet.Cells.Item(8, 3).Value = "长盛竞争优势股票C"
Set-TextValue $q4Sheet 8 4 "0.40"
Set-TextValue $q4Sheet 8 5 "88.91"
Set-TextValue $q4Sheet 8 6 "5.28"
Set-TextValue $q4Sheet 8 7 "0.0211"
$q4Sheet.Cells.Item(8, 8).Value = 6

# ---- row 9 ----
$q4Sheet.Cells.Item(9, 1).Value = 7
Set-TextValue $q4Sheet 9 2 "013881"
$q4Sheet.Cells.Item(9, 3).Value = "长信中证500指数增强C"
Set-TextValue $q4Sheet 9 4 "0.96"
Set-TextValue $q4Sheet 9 5 "92.79"
Set-TextValue $q4Sheet 9 6 "1.61"
Set-TextValue $q4Sheet 9 7 "0.0155"
$q4Sheet.Cells.Item(9, 8).Value = 6

# Leave the workbook selection the way it started (first/"总计" tab active).
$totalSheet.Activate()
$totalSheet.Range("A1").Select() | Out-Null
